$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: NW-621 entry
$ws.Range("A7").Value = "NW-621"
$ws.Range("B7").Value = "Network Protocol for 621"
$ws.Range("C7").Value = "-"
$ws.Range("D7").Value = "TCP/IP"

# Row 8: NW-622 entry
$ws.Range("A8").Value = "NW-622"
$ws.Range("B8").Value = "Network Protocol for 622"
$ws.Range("C8").Value = "-"
$ws.Range("D8").Value = "TCP/IP"

# D7/D8 reuse the same cell style as D6 (Arial font variant) - copy
# the existing style instead of creating a new one via Font assignment
$ws.Range("D6").Copy()
$null = $ws.Range("D7:D8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection to match the new extent of the table
$null = $ws.Range("A7:D8").Select()
